# Fix du téléchargement des toutes les sessions
#
# The campaign export used to dump a single, mis-shaped sample row per
# sheet. This rebuilds Sheet1 with the real per-session export shape: a
# bold / bordered header row (Token, Coords, Association Ligne,
# Association Col, Moyenne, Temps total, Choix Final, Id Campagne)
# followed by one row per recorded token/campaign, so every session in
# the upload actually gets written out instead of just the last one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (bold, thin border, centered / top-aligned) ---------
$ws.Cells.Item(1,1).Value = "Token"
$ws.Cells.Item(1,2).Value = "Coords"
$ws.Cells.Item(1,3).Value = "Association Ligne"
$ws.Cells.Item(1,4).Value = "Association Col"
$ws.Cells.Item(1,5).Value = "Moyenne"
$ws.Cells.Item(1,6).Value = "Temps total"
$ws.Cells.Item(1,7).Value = "Choix Final"
$ws.Cells.Item(1,8).Value = "Id Campagne"

$headerRange = $ws.Range("A1:H1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# ---- Row 2: campaign C3ID0efdcd --------------------------------------
$ws.Cells.Item(2,1).Value = "C3ID0efdcd"
$ws.Cells.Item(2,2).Value = "1:2', 2:3', 3:4', 3:3', 2:4'"
$ws.Cells.Item(2,3).Value = "Produit 4', Produit 3', Produit 1', Produit 3', Produit 1'"
$ws.Cells.Item(2,4).Value = "Lettre', Mot', Chiffre', Chiffre', Mot'"
$ws.Cells.Item(2,5).Value = 589
$ws.Cells.Item(2,6).Value = "Produit 1"
$ws.Cells.Item(2,7).Value = 3
$ws.Cells.Item(2,8).Value = "Token"

# ---- Row 3: campaign C3ID93a601 --------------------------------------
$ws.Cells.Item(3,1).Value = "C3ID93a601"
$ws.Cells.Item(3,2).Value = "2:1', 2:2', 2:3', 2:4', 1:4', 1:3'"
$ws.Cells.Item(3,3).Value = "Produit 3', Produit 2', Produit 1', Produit 4', Produit 4', Produit 1'"
$ws.Cells.Item(3,4).Value = "Lettre', Lettre', Lettre', Lettre', Chiffre', Chiffre'"
$ws.Cells.Item(3,5).Value = 418.3333333333333
$ws.Cells.Item(3,6).Value = "Produit 1"
$ws.Cells.Item(3,7).Value = 3
$ws.Cells.Item(3,8).Value = "Token"

# ---- Row 4: campaign C3ID651bab --------------------------------------
$ws.Cells.Item(4,1).Value = "C3ID651bab"
$ws.Cells.Item(4,2).Value = "1:1', 2:2', 2:3', 3:3', 3:4', 2:4', 1:4'"
$ws.Cells.Item(4,3).Value = "Produit 3', Produit 2', Produit 1', Produit 1', Produit 4', Produit 4', Produit 4'"
$ws.Cells.Item(4,4).Value = "Chiffre', Mot', Mot', Lettre', Lettre', Mot', Chiffre'"
$ws.Cells.Item(4,5).Value = 491.8571428571428
$ws.Cells.Item(4,6).Value = "Produit 4"
$ws.Cells.Item(4,7).Value = 3
$ws.Cells.Item(4,8).Value = "Token"

# ---- Row 5: carry-over of row 4's association data -------------------
$ws.Cells.Item(5,1).Value = "Token"
$ws.Cells.Item(5,2).Value = "1:1', 2:2', 2:3', 3:3', 3:4', 2:4', 1:4'"
$ws.Cells.Item(5,3).Value = "Produit 3', Produit 2', Produit 1', Produit 1', Produit 4', Produit 4', Produit 4'"
$ws.Cells.Item(5,4).Value = "Chiffre', Mot', Mot', Lettre', Lettre', Mot', Chiffre'"
$ws.Cells.Item(5,5).Value = "Token"
$ws.Cells.Item(5,6).Value = "Produit 4"
$ws.Cells.Item(5,7).Value = 3
$ws.Cells.Item(5,8).Value = "Token"

# ---- Row 6: carry-over of row 4's association data (repeat) ----------
$ws.Cells.Item(6,1).Value = "Token"
$ws.Cells.Item(6,2).Value = "1:1', 2:2', 2:3', 3:3', 3:4', 2:4', 1:4'"
$ws.Cells.Item(6,3).Value = "Produit 3', Produit 2', Produit 1', Produit 1', Produit 4', Produit 4', Produit 4'"
$ws.Cells.Item(6,4).Value = "Chiffre', Mot', Mot', Lettre', Lettre', Mot', Chiffre'"
$ws.Cells.Item(6,5).Value = "Token"
$ws.Cells.Item(6,6).Value = "Produit 4"
$ws.Cells.Item(6,7).Value = 3
$ws.Cells.Item(6,8).Value = "Token"

# ---- Rows 7-9: trailing blank rows kept in the export range -----------
for ($r = 7; $r -le 9; $r++) {
    for ($c = 1; $c -le 8; $c++) {
        $ws.Cells.Item($r, $c).Value = "Token"
    }
}
